# daily auto push: 2026-01-29 04:13 UTC
# Insert a new data row at row 744 ("2026/01/29", 木, 12, 159), shifting all
# subsequent rows (old 744..785) down by one (new 745..786).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 744, pushing everything below down.
$ws.Rows.Item(744).Insert()

# Fill in the new row's data.
# The date is written with a leading apostrophe to force text entry (avoids
# the host auto-converting the "yyyy/mm/dd" literal into a date serial
# number), then formatting is cleared so the cell keeps the default style
# (matching the unstyled data cells elsewhere in the column).
$ws.Range("A744").Value = "'2026/01/29"
$ws.Range("A744").ClearFormats()

$ws.Range("B744").Value = "木"
$ws.Range("C744").Value = 12
$ws.Range("D744").Value = 159
